# Applies the cryptos-list price/volume refresh described in the commit
# message ("Updated cryptos list ... with GitHub Actions"). Each cell is
# forced to Text format before the write (and the style is reset to Normal
# right after) so numeric-looking strings such as "332.61" are stored as
# text -- exactly like the original inline-string cells -- instead of being
# auto-coerced into numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = '27.587.44'
    "E2" = '  -1.48%  '
    "D3" = '1.847.89'
    "E3" = '  -2.27%  '
    "E4" = '  -1.18%  '
    "D5" = '332.61'
    "E5" = '  -1.00%  '
    "E6" = '  -1.31%  '
    "D7" = '0.4641'
    "E7" = '  -1.47%  '
    "D8" = '0.3854'
    "E8" = '  -2.37%  '
    "D9" = '46.07'
    "E9" = '  -1.54%  '
    "D10" = '0.07916'
    "E10" = '  -1.24%  '
    "D11" = '0.9949'
    "E11" = '  -2.17%  '
    "D12" = '21.47'
    "E12" = '  -1.54%  '
    "D13" = '1.838.36'
    "E13" = '  -3.02%  '
    "D14" = '5.921'
    "E14" = '  -1.11%  '
    "D16" = '1.003'
    "E16" = '  -1.48%  '
    "D17" = '88.77'
    "E17" = '  +0.87%  '
    "D18" = '0.06653'
    "E18" = '  -2.02%  '
    "E19" = '  -1.55%  '
    "D20" = '17.06'
    "E20" = '  -0.76%  '
    "E21" = '  -1.23%  '
    "D22" = '27.570.70'
    "E22" = '  -1.58%  '
    "D23" = '5.385'
    "E23" = '  -2.09%  '
    "D24" = '10.92'
    "E24" = '  -0.60%  '
    "D25" = '2.298'
    "E25" = '  -2.75%  '
    "D26" = '2.065.51'
    "E26" = '  -2.27%  '
    "D27" = '157.87'
    "E27" = '  -0.96%  '
    "D28" = '19.52'
    "E28" = '  -2.52%  '
    "D29" = '2.103'
    "E29" = '  +0.01%  '
    "D30" = '5.407'
    "E30" = '  -1.54%  '
    "D31" = '119.78'
    "E31" = '  -1.50%  '
    "D32" = '0.9754'
    "E32" = '  +1.01%  '
    "D33" = '0.09397'
    "E33" = '  -1.78%  '
    "D34" = '3.576'
    "E34" = '  -1.96%  '
    "D35" = '5.286'
    "E35" = '  -1.30%  '
    "D36" = '1.347'
    "E36" = '  -1.14%  '
    "D37" = '0.06012'
    "E37" = '  -1.83%  '
    "D38" = '0.02221'
    "E38" = '  -1.20%  '
    "D39" = '8.288'
    "E39" = '  +0.67%  '
    "D40" = '1.182'
    "E40" = '  -2.82%  '
    "D41" = '0.5901'
    "E41" = '  -1.04%  '
    "D42" = '0.1862'
    "E42" = '  -2.15%  '
    "D43" = '10.31'
    "E43" = '  -0.24%  '
    "D44" = '1.246'
    "E44" = '  -2.01%  '
    "D45" = '0.5584'
    "E45" = '  -2.05%  '
    "D46" = '12.11'
    "E46" = '  -0.95%  '
    "D47" = '1.897'
    "E47" = '  -2.47%  '
    "E48" = '  -2.58%  '
    "D49" = '110.66'
    "E49" = '  -2.56%  '
    "D50" = '1.053'
    "E50" = '  -1.65%  '
    "D51" = '1.001'
    "E51" = '  -1.46%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
